$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "59.107.56"
$ws.Range("E2").Value = "  -0.29%  "
$ws.Range("D3").Value = "2.644.49"
$ws.Range("E3").Value = "  +0.23%  "
$ws.Range("E4").Value = "  +0.04%  "
$ws.Range("D5").Value = "'522.32"
$ws.Range("E5").Value = "  +1.34%  "
$ws.Range("D6").Value = "'145.35"
$ws.Range("E6").Value = "  -0.20%  "
$ws.Range("D7").Value = "'1.00"
$ws.Range("E7").Value = "  +0.08%  "
$ws.Range("E8").Value = "  -0.19%  "
$ws.Range("D9").Value = "2.649.20"
$ws.Range("E9").Value = "  +0.14%  "
$ws.Range("D10").Value = "'6.92"
$ws.Range("E10").Value = "  +9.91%  "
$ws.Range("E11").Value = "  -2.38%  "
$ws.Range("E12").Value = "  -0.31%  "
$ws.Range("E13").Value = "  +2.05%  "
$ws.Range("D14").Value = "3.110.38"
$ws.Range("E14").Value = "  +0.09%  "
$ws.Range("D15").Value = "59.631.84"
$ws.Range("E15").Value = "  +0.51%  "
$ws.Range("D16").Value = "'21.16"
$ws.Range("E16").Value = "  +0.88%  "
$ws.Range("B17").Value = "ShibaInu"
$ws.Range("C17").Value = "https://coinranking.com/coin/xz24e0BjL+shibainu-shib"
$ws.Range("D17").Value = "'0.0000136"
$ws.Range("E17").Value = "  -0.67%  "
$ws.Range("B18").Value = "WrappedEther"
$ws.Range("C18").Value = "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
$ws.Range("D18").Value = "2.647.34"
$ws.Range("E18").Value = "  -0.16%  "
$ws.Range("D19").Value = "'342.23"
$ws.Range("E19").Value = "  -1.97%  "
$ws.Range("E20").Value = "  -1.47%  "
$ws.Range("D21").Value = "'10.36"
$ws.Range("E21").Value = "  +0.50%  "
$ws.Range("D22").Value = "'6.31"
$ws.Range("E22").Value = "  +1.47%  "
$ws.Range("D23").Value = "'1.00"
$ws.Range("E23").Value = "  +0.08%  "
$ws.Range("D24").Value = "'63.65"
$ws.Range("E24").Value = "  +2.30%  "
$ws.Range("E25").Value = "  +1.13%  "
$ws.Range("E26").Value = "  -0.68%  "
$ws.Range("E27").Value = "  +0.38%  "
$ws.Range("D28").Value = "0.0₃0804"
$ws.Range("E28").Value = "  +0.32%  "
$ws.Range("D29").Value = "'7.10"
$ws.Range("E29").Value = "  -0.65%  "
$ws.Range("D30").Value = "'6.68"
$ws.Range("E30").Value = "  +3.29%  "
$ws.Range("E31").Value = "  -0.10%  "
$ws.Range("D32").Value = "'1.59"
$ws.Range("E32").Value = "  +0.85%  "
$ws.Range("D33").Value = "'18.77"
$ws.Range("D34").Value = "'149.33"
$ws.Range("E34").Value = "  -0.33%  "
$ws.Range("D35").Value = "'4.19"
$ws.Range("E35").Value = "  +3.25%  "
$ws.Range("E36").Value = "  +2.32%  "
$ws.Range("D37").Value = "'0.898"
$ws.Range("E37").Value = "  -5.19%  "
$ws.Range("D38").Value = "'0.883"
$ws.Range("E38").Value = "  +2.75%  "
$ws.Range("D39").Value = "'36.76"
$ws.Range("E39").Value = "  +0.39%  "
$ws.Range("D40").Value = "'1.48"
$ws.Range("E40").Value = "  +0.47%  "
$ws.Range("E41").Value = "  -1.95%  "
$ws.Range("E42").Value = "  +4.59%  "
$ws.Range("E43").Value = "  +0.21%  "
$ws.Range("D44").Value = "'275.23"
$ws.Range("E44").Value = "  -0.90%  "
$ws.Range("E45").Value = "  +1.02%  "
$ws.Range("E46").Value = "  -1.59%  "
$ws.Range("D47").Value = "'0.0539"
$ws.Range("E47").Value = "  +1.80%  "
$ws.Range("D48").Value = "2.055.21"
$ws.Range("E48").Value = "  -0.90%  "
$ws.Range("E49").Value = "  +1.97%  "
$ws.Range("D50").Value = "'4.78"
$ws.Range("E50").Value = "  +1.25%  "
$ws.Range("D51").Value = "'19.08"
$ws.Range("E51").Value = "  +0.86%  "
